$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new backlog items to the task list (rows 19 and 20)
$ws.Range("A19").Value = "Gamepad Input"
$ws.Range("B19").Value = 3

$ws.Range("A20").Value = "Fix Dialogue"
$ws.Range("B20").Value = 6

# Update the selected cell to match the author's saved selection
$ws.Range("A3").Select()
